$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 283; this pushes the existing rows
# 283-358 down to 284-359 (growing the used range to A1:R359).
$ws.Rows.Item(283).Insert()

# Populate the newly inserted row 283 with its data.
$ws.Cells.Item(283, 1).Value = 3
$ws.Cells.Item(283, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(283, 3).Value = "Coquimbo"
$ws.Cells.Item(283, 4).Value = 44642
$ws.Cells.Item(283, 5).Value = 5
$ws.Cells.Item(283, 6).Value = 100112017
$ws.Cells.Item(283, 7).Value = "Apio"
$ws.Cells.Item(283, 8).Value = "Americana (o)"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 190
$ws.Cells.Item(283, 11).Value = 9500
$ws.Cells.Item(283, 12).Value = 10000
$ws.Cells.Item(283, 13).Value = 9789
$ws.Cells.Item(283, 14).Value = "`$/docena de matas"
$ws.Cells.Item(283, 15).Value = "Pan de Az$([char]0x00FA)car"
$ws.Cells.Item(283, 16).Value = 1632
$ws.Cells.Item(283, 17).Value = 6
$ws.Cells.Item(283, 18).Value = "Hortaliza"
